$d = $word.ActiveDocument

# Correção ortográfica: no item "P2. O usuário define quais serão as
# pesquisas rápidas;" o ";" final deve virar ".".
# Esse é o único ";" do documento.
$rng = $d.Content
$found = $rng.Find.Execute(";", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    # O Word move automaticamente a marca oculta "_GoBack" para o ponto da
    # última edição feita no documento. Reproduzimos isso marcando a posição
    # do ";" (que será o local da edição) antes de trocar o texto - assim a
    # marca "_GoBack" sai de onde estava (fim do documento) e passa a ficar
    # ao lado do texto corrigido.
    $d.Bookmarks.Add("_GoBack", $rng)

    $rng2 = $d.Content
    $rng2.Find.Execute(";", $true, $false, $false, $false, $false, `
                        $true, 1, $false, ".", 2)
}
